$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 237
$ws.Range("A237").Value = "What's the maximum curves I can load in a data file?"
$ws.Range("B237").Value = "llama3.2:latest"
$ws.Range("C237").Value = "According to Document 30, the maximum number of data files you can load is unlimited. However, it does not specify the maximum number of curves per data file.`nTo answer your question accurately, let's look at another document that might provide more information on this topic.`nDocument 29 states: ""Deviation from the above may cause errors during loading data from files."" It also mentions that if possible, you should ask for LAS files from your data provider. However, it does not specify a maximum number of curves per file.`nUnfortunately, without further information or clarification, I cannot provide an accurate answer to your question about the maximum number of curves you can load in a single data file."
$ws.Rows.Item(237).AutoFit()

# Row 238
$ws.Range("A238").Value = "Why can I not add another layout to my log?"
$ws.Range("B238").Value = "llama3.2:latest"
$ws.Range("C238").Value = "According to the document, if adding a layout to an existing 'blank' session of GEO, you should load an existing Vew file or create a new one. If adding a layout to an odf with layouts already defined, you need to create the layout first (step 1-2), which will create a copy of the layout currently opened, and then you can overwrite its contents by loading a view file or creating a new one by editing the existing tracks."

# Row 239
$ws.Range("A239").Value = "Why can I not add another layout to my log?"
$ws.Range("B239").Value = "llama3.2:latest"
$ws.Range("C239").Value = "According to the document, if adding a layout to an existing 'blank' session of GEO, you should load an existing Vew file or create a new one. If adding a layout to an odf with layouts already defined, you need to create the layout first (step 1-2), which will create a copy of the layout currently opened, and then you can overwrite its contents by loading a view file or creating a new one by editing the existing tracks."

# Row 240
$ws.Range("A240").Value = "Why can I not add another layout to my log?"
$ws.Range("B240").Value = "llama3.2:latest"
$ws.Range("C240").Value = "According to the document, you cannot add another layout to your log because the layout specifies up to 19 layouts per ODF file."

# Row 241
$ws.Range("A241").Value = "Why can't I add 251 curve shades to my log?"
$ws.Range("B241").Value = "llama3.2:latest"
$ws.Range("C241").Value = "According to document 25, the maximum number of curve shades per plot is 250. This means you cannot add more than 250 curve shades to your log."

# Row 242
$ws.Range("A242").Value = "I have 20000 modifiers added ty log, why I can't I add anymore?"
$ws.Range("B242").Value = "llama3.2:latest"
$ws.Range("C242").Value = "You cannot add more than 20000 modifiers per plot because of the limit specified on theHometab."
